# Speech.docx plan update:
# Replace the trailing "4." / "5." placeholder paragraphs with the new
# numbered plan items, and append a trailing blank paragraph.

$d = $word.ActiveDocument

# Locate the two placeholder paragraphs ("4." and "5.") - they are the
# last two paragraphs of the document body, right before the sectPr.
$count = $d.Paragraphs.Count
$p1 = $d.Paragraphs.Item($count - 1)
$p2 = $d.Paragraphs.Item($count)

$start = $p1.Range.Start
$end = $p2.Range.End
$target = $d.Range($start, $end)

$W = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$xml = @'
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr>
    <w:rPr><w:lang w:val="ru-RU"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>Замикання та область видимості.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr>
    <w:rPr><w:lang w:val="ru-RU"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>Об'єкти і методи.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr>
    <w:rPr><w:lang w:val="ru-RU"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>Методи в об'єктів. Перетворення об'єктів.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr>
    <w:rPr><w:lang w:val="ru-RU"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">Створення об'єктів через </w:t></w:r>
  <w:r><w:rPr><w:i/></w:rPr><w:t>new</w:t></w:r>
  <w:r><w:rPr><w:i/><w:lang w:val="uk-UA"/></w:rPr><w:t>.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr>
    <w:rPr><w:lang w:val="ru-RU"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>Дескриптори, геттери і сеттери.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr>
    <w:rPr><w:lang w:val="ru-RU"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:lastRenderedPageBreak/><w:t>Статичні методи. Фабричні методи.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr>
    <w:rPr><w:lang w:val="ru-RU"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="uk-UA"/></w:rPr><w:t xml:space="preserve">Виклики методів: </w:t></w:r>
  <w:r><w:rPr><w:i/></w:rPr><w:t>call, apply.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr>
    <w:rPr><w:lang w:val="ru-RU"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="uk-UA"/></w:rPr><w:t>Прив'язка контексту</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr>
    <w:rPr><w:lang w:val="ru-RU"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="uk-UA"/></w:rPr><w:t>Декоратори.</w:t></w:r>
  <w:bookmarkStart w:id="2" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="2"/>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:left="720" w:firstLine="0"/>
    <w:rPr><w:lang w:val="uk-UA"/></w:rPr>
  </w:pPr>
</w:p>
'@

# The curly apostrophe (U+2019) in "Прив'язка" must be the typographic
# quote used in the source document, not the straight ASCII apostrophe.
$xml = $xml.Replace("Прив'язка", "Прив" + [char]0x2019 + "язка")

$target.InsertXML($xml)

# Update the ToC page reference for "ЗМІСТ" (the first PAGEREF field)
# from 1 to 3, reflecting the document now spanning more pages.
$tocField = $d.Fields.Item(2)
Write-Host "Field2 code: $($tocField.Code.Text)"
